$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# New numeric data used as the charts' source range (Folha1!B1:B4 series,
# plus a couple of stray values elsewhere on the sheet).
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 1
$ws1.Range("A4").Value = 1
$ws1.Range("A5").Value = 1
$ws1.Range("B8").Value = 1
$ws1.Range("D10").Value = 1

# Multi-line note in A1, wrapped, with a taller row to match.
$ws1.Range("A1").Value = "asdasdasdsda" + [char]13 + [char]10 + "asd" + [char]13 + [char]10 + "asd" + [char]13 + [char]10
$ws1.Range("A1").WrapText = $true
$ws1.Rows(1).RowHeight = 75

# Page setup for this sheet.
$ws1.PageSetup.PaperSize = 9
$ws1.PageSetup.Orientation = 1

# Two line charts plotting Folha1!$B$2:$B$4 with the series name from $B$1.
$chart2 = $ws1.ChartObjects().Add(338.75, 76.25, 300, 300)
$chart2.Name = "Grafico 2"
$chart2.Chart.ChartType = 4
$chart2.Chart.SetSourceData($ws1.Range("B1:B4"))
$chart2.Chart.HasTitle = $true
$chart2.Chart.ChartTitle.Text = "Graph Title"
$chart2.Chart.HasLegend = $true
$chart2.Chart.Legend.Position = -4152
$chart2.Chart.Axes(1).HasTitle = $true
$chart2.Chart.Axes(1).AxisTitle.Text = "Title of X axis"
$chart2.Chart.Axes(2).HasTitle = $true
$chart2.Chart.Axes(2).AxisTitle.Text = "Title of Y axis"

$chart1 = $ws1.ChartObjects().Add(686.75, 90.5, 300, 300)
$chart1.Name = "Grafico 1"
$chart1.Chart.ChartType = 4
$chart1.Chart.SetSourceData($ws1.Range("B1:B4"))
$chart1.Chart.HasTitle = $true
$chart1.Chart.ChartTitle.Text = "Graph Title"
$chart1.Chart.HasLegend = $true
$chart1.Chart.Legend.Position = -4152
$chart1.Chart.Axes(1).HasTitle = $true
$chart1.Chart.Axes(1).AxisTitle.Text = "Title of X axis"
$chart1.Chart.Axes(2).HasTitle = $true
$chart1.Chart.Axes(2).AxisTitle.Text = "Title of Y axis"

# Select D10 then make Folha1 the active sheet/tab (it was Folha5 before).
$ws1.Range("D10").Select()
$ws1.Activate()

Write-Output "done"
